# Auto-generated edit script for 江西-漫展信息.xlsx
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("A2").Value = 1
    $ws.Range("B2").Value = "2024.02.24"
    $ws.Range("C2").Value = "宜春·融荟城难忘今宵汉文化节"
    $ws.Range("D2").Value = "宜阳大道239号 宜春融荟城"
    $ws.Range("E2").Value = "2024.02.24 14:00-02.24 18:00"
    $ws.Range("F2").Value = 29
    $ws.Range("G2").Value = "已停售"
    $ws.Range("H2").Value = "https://show.bilibili.com/platform/detail.html?id=81690"
    $ws.Range("I2").Value = "//i0.hdslb.com/bfs/openplatform/202402/ldtkc9Sp1706865634128.jpeg"

    $ws.Range("A3").Value = 2
    $ws.Range("B3").Value = "2024.02.24"
    $ws.Range("C3").Value = "景德镇·陶溪川×次元文化元宵游园会（ 免费活动）"
    $ws.Range("D3").Value = "新厂西路315号 陶溪川发布大厅"
    $ws.Range("E3").Value = "2024.02.24 10:00-02.25 18:00"
    $ws.Range("F3").Value = 429
    $ws.Range("G3").Value = 30
    $ws.Range("H3").Value = "https://show.bilibili.com/platform/detail.html?id=81207"
    $ws.Range("I3").Value = "//i1.hdslb.com/bfs/openplatform/202402/nIs2jtUn1707298876430.png"

    $ws.Range("A4").Value = 3
    $ws.Range("B4").Value = "2024.03.02"
    $ws.Range("C4").Value = "南昌·meeting动漫游戏嘉年华"
    $ws.Range("D4").Value = "南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆"
    $ws.Range("E4").Value = "2024.03.02 09:00-03.03 17:00"
    $ws.Range("F4").Value = 1447
    $ws.Range("G4").Value = 60
    $ws.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=79555"
    $ws.Range("I4").Value = "//i0.hdslb.com/bfs/openplatform/202402/l6GUtggC1706843695971.jpeg"

    $ws.Range("A5").Value = 4
    $ws.Range("B5").Value = "2024.03.09"
    $ws.Range("C5").Value = "景德镇·江报国风动漫展 "
    $ws.Range("D5").Value = "迎宾大道与寺山路交叉口东200米 陶博城"
    $ws.Range("E5").Value = "2024.03.09 09:00-03.10 17:00"
    $ws.Range("F5").Value = 945
    $ws.Range("G5").Value = 55
    $ws.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=81362"
    $ws.Range("I5").Value = "//i2.hdslb.com/bfs/openplatform/202402/oM49o66R1708334630235.jpeg"

    $ws.Range("A6").Value = 5
    $ws.Range("B6").Value = "2024.03.16"
    $ws.Range("C6").Value = "景德镇·原神X崩铁X崩坏动漫展only"
    $ws.Range("D6").Value = "陶阳南路188号 晨枫臻品酒店"
    $ws.Range("E6").Value = "2024.03.16 10:00-03.16 17:00"
    $ws.Range("F6").Value = 60
    $ws.Range("G6").Value = 55
    $ws.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=80920"
    $ws.Range("I6").Value = "//i0.hdslb.com/bfs/openplatform/202401/IugBckTp1705469476482.png"

    $ws.Range("A7").Value = 6
    $ws.Range("B7").Value = "2024.03.16"
    $ws.Range("C7").Value = "江西·ShiningStaR动漫游戏文化节5th"
    $ws.Range("D7").Value = "高新开发区紫阳大道666号 江西奥林匹克体育中心综合训练馆"
    $ws.Range("E7").Value = "2024.03.16 09:30-03.17 17:00"
    $ws.Range("F7").Value = 2074
    $ws.Range("G7").Value = 60
    $ws.Range("H7").Value = "https://show.bilibili.com/platform/detail.html?id=81792"
    $ws.Range("I7").Value = "//i2.hdslb.com/bfs/openplatform/202402/2l16aHBJ1707209383729.jpeg"

    $ws.Range("A8").Value = 7
    $ws.Range("B8").Value = "2024.03.23"
    $ws.Range("C8").Value = "上饶·原×铁×崩only"
    $ws.Range("D8").Value = "五三东大道42号 回禾酒店"
    $ws.Range("E8").Value = "2024.03.23 10:00-03.23 17:00"
    $ws.Range("F8").Value = 34
    $ws.Range("G8").Value = 60
    $ws.Range("H8").Value = "https://show.bilibili.com/platform/detail.html?id=81103"
    $ws.Range("I8").Value = "//i2.hdslb.com/bfs/openplatform/202401/pp6c5TsC1705647180602.jpeg"

    $ws.Range("A9").Value = 8
    $ws.Range("B9").Value = "2024.03.23"
    $ws.Range("C9").Value = "南昌·AP动漫游戏嘉年华"
    $ws.Range("D9").Value = "八一桥街道青山南路118号 蓝海会展中心"
    $ws.Range("E9").Value = "2024.03.23 09:00-03.24 17:00"
    $ws.Range("F9").Value = 1256
    $ws.Range("G9").Value = 58.5
    $ws.Range("H9").Value = "https://show.bilibili.com/platform/detail.html?id=81232"
    $ws.Range("I9").Value = "//i2.hdslb.com/bfs/openplatform/202401/NZv97SmS1705912230957.jpeg"

    $ws.Range("A10").Value = 9
    $ws.Range("B10").Value = "2024.03.23"
    $ws.Range("C10").Value = "南昌·原X穹X崩only"
    $ws.Range("D10").Value = "丰和北大道299号 新吉花园酒店"
    $ws.Range("E10").Value = "2024.03.23 10:00-03.23 17:00"
    $ws.Range("F10").Value = 62
    $ws.Range("G10").Value = 65
    $ws.Range("H10").Value = "https://show.bilibili.com/platform/detail.html?id=80807"
    $ws.Range("I10").Value = "//i0.hdslb.com/bfs/openplatform/202401/rY4v2Opx1705051458246.jpeg"

    $ws.Range("A11").Value = 10
    $ws.Range("B11").Value = "2024.03.23"
    $ws.Range("C11").Value = "南昌·运动番only春季集训"
    $ws.Range("D11").Value = "创新三路777号 南昌小飞侠章鱼文化体育公园"
    $ws.Range("E11").Value = "2024.03.23 10:00-03.24 17:00"
    $ws.Range("F11").Value = 103
    $ws.Range("G11").Value = 58
    $ws.Range("H11").Value = "https://show.bilibili.com/platform/detail.html?id=81950"
    $ws.Range("I11").Value = "//i1.hdslb.com/bfs/openplatform/202402/bm4uH4qB1708425538357.jpeg"

    $ws.Range("A12").Value = 11
    $ws.Range("B12").Value = "2024.03.24"
    $ws.Range("C12").Value = "南昌·AP动漫游戏  嘉年华内场票-小N&子音"
    $ws.Range("D12").Value = "八一桥街道青山南路118号 蓝海会展中心"
    $ws.Range("E12").Value = "2024.03.24 09:00-03.24 17:00"
    $ws.Range("F12").Value = 35
    $ws.Range("G12").Value = 218
    $ws.Range("H12").Value = "https://show.bilibili.com/platform/detail.html?id=81973"
    $ws.Range("I12").Value = "//i0.hdslb.com/bfs/openplatform/202402/zbG5HICL1708504962467.jpeg"

    $ws.Range("A13").Value = 12
    $ws.Range("B13").Value = "2024.03.30"
    $ws.Range("C13").Value = "南昌·CM01动漫游戏博览会"
    $ws.Range("D13").Value = "怀玉山大道1315号 南昌绿地国际博览中心"
    $ws.Range("E13").Value = "2024.03.30 10:00-03.31 17:00"
    $ws.Range("F13").Value = 303
    $ws.Range("G13").Value = 55
    $ws.Range("H13").Value = "https://show.bilibili.com/platform/detail.html?id=81691"
    $ws.Range("I13").Value = "//i2.hdslb.com/bfs/openplatform/202402/IYLaH7AS1706866218597.png"

    $ws.Range("A14").Value = 13
    $ws.Range("B14").Value = "2024.03.30"
    $ws.Range("C14").Value = "鹰潭·原×铁×崩only"
    $ws.Range("D14").Value = "南站路24号 回禾酒店(鹰潭火车站店)"
    $ws.Range("E14").Value = "2024.03.30 10:00-03.30 17:00"
    $ws.Range("F14").Value = 25
    $ws.Range("G14").Value = 60
    $ws.Range("H14").Value = "https://show.bilibili.com/platform/detail.html?id=81097"
    $ws.Range("I14").Value = "//i2.hdslb.com/bfs/openplatform/202401/q0AZaXAk1705646244207.jpeg"

    # Remove the now-obsolete trailing rows (old rows 15 and 16)
    $ws.Rows.Item(15).Delete() | Out-Null
    $ws.Rows.Item(15).Delete() | Out-Null
}

$wb.Save()